$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source sheet is a Google-Forms/Sheets export where each row carries an
# explicit banded-row style (alternating cellXfs per row parity, with a unique
# border style reserved for whatever is currently the last row of the table).
# New rows 268-293 are appended, row 267 stops being the last row (and loses
# its special border), and row 293 becomes the new last row and inherits that
# special border instead. We replicate the styling by copying whole rows from
# existing rows that already carry the desired style, then overwrite the values.

# 1) Snapshot row 267 (current last row, special border style) into the new last
#    row 293 before row 267 itself gets restyled.
$ws.Range("A267:N267").Copy($ws.Range("A293:N293"))
$ws.Range("L267").Copy($ws.Range("M293"))
$ws.Range("N293").Clear()

# 2) Restyle row 267 and create rows 268-292 by copying the matching banded
#    template row already present on the sheet (rows 4, 5, 7 and 10 cover the
#    four style/column combinations used throughout the table).
$ws.Range("A5:N5").Copy($ws.Range("A267:N267"))
$ws.Range("A10:N10").Copy($ws.Range("A268:N268"))
$ws.Range("A7:N7").Copy($ws.Range("A269:N269"))
$ws.Range("A4:N4").Copy($ws.Range("A270:N270"))
$ws.Range("A5:N5").Copy($ws.Range("A271:N271"))
$ws.Range("A4:N4").Copy($ws.Range("A272:N272"))
$ws.Range("A7:N7").Copy($ws.Range("A273:N273"))
$ws.Range("A4:N4").Copy($ws.Range("A274:N274"))
$ws.Range("A5:N5").Copy($ws.Range("A275:N275"))
$ws.Range("A10:N10").Copy($ws.Range("A276:N276"))
$ws.Range("A7:N7").Copy($ws.Range("A277:N277"))
$ws.Range("A10:N10").Copy($ws.Range("A278:N278"))
$ws.Range("A7:N7").Copy($ws.Range("A279:N279"))
$ws.Range("A10:N10").Copy($ws.Range("A280:N280"))
$ws.Range("A5:N5").Copy($ws.Range("A281:N281"))
$ws.Range("A4:N4").Copy($ws.Range("A282:N282"))
$ws.Range("A7:N7").Copy($ws.Range("A283:N283"))
$ws.Range("A10:N10").Copy($ws.Range("A284:N284"))
$ws.Range("A7:N7").Copy($ws.Range("A285:N285"))
$ws.Range("A10:N10").Copy($ws.Range("A286:N286"))
$ws.Range("A5:N5").Copy($ws.Range("A287:N287"))
$ws.Range("A4:N4").Copy($ws.Range("A288:N288"))
$ws.Range("A5:N5").Copy($ws.Range("A289:N289"))
$ws.Range("A10:N10").Copy($ws.Range("A290:N290"))
$ws.Range("A5:N5").Copy($ws.Range("A291:N291"))
$ws.Range("A10:N10").Copy($ws.Range("A292:N292"))

# 3) Write the real values/text into rows 267 (unchanged data, restyled only)
#    through 293 (brand-new survey responses).
# row 267
$ws.Range("A267").Value = 45570.45796989583
$ws.Range("B267").Value = 'mt1661@naver.com'
$ws.Range("C267").Value = '콘탠츠 IT'
$ws.Range("D267").Value = 20215239
$ws.Range("E267").Value = '정성민'
$ws.Range("F267").Value = '민주 문자'
$ws.Range("G267").Value = '한자'
$ws.Range("H267").Value = '3개'
$ws.Range("I267").Value = 0.8
$ws.Range("J267").Value = '대한민국'
$ws.Range("K267").Value = '2배 정도 실직할 가능성이 높다'
$ws.Range("L267").Value = 'Black'
$ws.Range("N267").Value = '모름/기타'
# row 268
$ws.Range("A268").Value = 45570.48767313657
$ws.Range("B268").Value = 'a35142191@gmail.com'
$ws.Range("C268").Value = '데이터사이언스학부'
$ws.Range("D268").Value = 20243241
$ws.Range("E268").Value = '이윤재'
$ws.Range("F268").Value = '민주 문자'
$ws.Range("G268").Value = '한글'
$ws.Range("H268").Value = '1개'
$ws.Range("I268").Value = 0.5
$ws.Range("J268").Value = '대한민국'
$ws.Range("K268").Value = '남들을 덜 신뢰한다'
$ws.Range("L268").Value = 'Black'
$ws.Range("N268").Value = '헐, 반 밖에 안 남았네.'
# row 269
$ws.Range("A269").Value = 45570.49257619213
$ws.Range("B269").Value = 'wonda0322@naver.com'
$ws.Range("C269").Value = '글로벌학부'
$ws.Range("D269").Value = 20246414
$ws.Range("E269").Value = '원다연'
$ws.Range("F269").Value = '민주 문자'
$ws.Range("G269").Value = '한글'
$ws.Range("H269").Value = '1개'
$ws.Range("I269").Value = 0.8
$ws.Range("J269").Value = '대한민국'
$ws.Range("K269").Value = '시간당 중위 임금이 60% 낮다'
$ws.Range("L269").Value = 'Red'
$ws.Range("M269").Value = '헐, 반 밖에 안 남았네.'
# row 270
$ws.Range("A270").Value = 45570.51792457176
$ws.Range("B270").Value = 'thwls5541@gmail.com'
$ws.Range("C270").Value = '법학과'
$ws.Range("D270").Value = 20207065
$ws.Range("E270").Value = '김소원'
$ws.Range("F270").Value = '민주 문자'
$ws.Range("G270").Value = '한자'
$ws.Range("H270").Value = '1개'
$ws.Range("I270").Value = 0.5
$ws.Range("J270").Value = '이탈리아'
$ws.Range("K270").Value = '건강이 좋지 않다'
$ws.Range("L270").Value = 'Red'
$ws.Range("M270").Value = '휴우, 그래도 반이나 남았네.'
# row 271
$ws.Range("A271").Value = 45570.55450297454
$ws.Range("B271").Value = 'kyj57980@gmail.com'
$ws.Range("C271").Value = '사회복지학과'
$ws.Range("D271").Value = 20202319
$ws.Range("E271").Value = '김예진'
$ws.Range("F271").Value = '엘리트 문자'
$ws.Range("G271").Value = '한자'
$ws.Range("H271").Value = '하나도 없다'
$ws.Range("I271").Value = 0.5
$ws.Range("J271").Value = '미국'
$ws.Range("K271").Value = '사회활동이나 자원활동에 덜 참여한다'
$ws.Range("L271").Value = 'Black'
$ws.Range("N271").Value = '헐, 반 밖에 안 남았네.'
# row 272
$ws.Range("A272").Value = 45570.55773021991
$ws.Range("B272").Value = 'dohan5453@gmail.com'
$ws.Range("C272").Value = '소프트웨어학부'
$ws.Range("D272").Value = 20225117
$ws.Range("E272").Value = '김도한'
$ws.Range("F272").Value = '민주 문자'
$ws.Range("G272").Value = '한글'
$ws.Range("H272").Value = '하나도 없다'
$ws.Range("I272").Value = 0.8
$ws.Range("J272").Value = '대한민국'
$ws.Range("K272").Value = '건강이 좋지 않다'
$ws.Range("L272").Value = 'Red'
$ws.Range("M272").Value = '휴우, 그래도 반이나 남았네.'
# row 273
$ws.Range("A273").Value = 45570.5690947338
$ws.Range("B273").Value = 'ahrang1225@gmail.com'
$ws.Range("C273").Value = '사회학과'
$ws.Range("D273").Value = 20242240
$ws.Range("E273").Value = '탁아랑'
$ws.Range("F273").Value = '민주 문자'
$ws.Range("G273").Value = '한글'
$ws.Range("H273").Value = '2개'
$ws.Range("I273").Value = 0.2
$ws.Range("J273").Value = '대한민국'
$ws.Range("K273").Value = '남들을 덜 신뢰한다'
$ws.Range("L273").Value = 'Red'
$ws.Range("M273").Value = '휴우, 그래도 반이나 남았네.'
# row 274
$ws.Range("A274").Value = 45570.581141180555
$ws.Range("B274").Value = 'ryul1128@naver.com'
$ws.Range("C274").Value = '식품영양학과'
$ws.Range("D274").Value = 20243806
$ws.Range("E274").Value = '김률아'
$ws.Range("F274").Value = '민주 문자'
$ws.Range("G274").Value = '한글'
$ws.Range("H274").Value = '하나도 없다'
$ws.Range("I274").Value = 0.8
$ws.Range("J274").Value = '대한민국'
$ws.Range("K274").Value = '시간당 중위 임금이 60% 낮다'
$ws.Range("L274").Value = 'Red'
$ws.Range("M274").Value = '휴우, 그래도 반이나 남았네.'
# row 275
$ws.Range("A275").Value = 45570.58408630787
$ws.Range("B275").Value = '123plokml123@gmail.com'
$ws.Range("C275").Value = '소프트웨어학부'
$ws.Range("D275").Value = 20245266
$ws.Range("E275").Value = '지현배'
$ws.Range("F275").Value = '민주 문자'
$ws.Range("G275").Value = '한글'
$ws.Range("H275").Value = '1개'
$ws.Range("I275").Value = 0.8
$ws.Range("J275").Value = '대한민국'
$ws.Range("K275").Value = '남들을 덜 신뢰한다'
$ws.Range("L275").Value = 'Black'
$ws.Range("N275").Value = '헐, 반 밖에 안 남았네.'
# row 276
$ws.Range("A276").Value = 45570.59573396991
$ws.Range("B276").Value = 'leeyubin050328@naver.com'
$ws.Range("C276").Value = '경영대학'
$ws.Range("D276").Value = 20243008
$ws.Range("E276").Value = '이유빈'
$ws.Range("F276").Value = '민주 문자'
$ws.Range("G276").Value = '한글'
$ws.Range("H276").Value = '하나도 없다'
$ws.Range("I276").Value = 0.8
$ws.Range("J276").Value = '대한민국'
$ws.Range("K276").Value = '남들을 덜 신뢰한다'
$ws.Range("L276").Value = 'Black'
$ws.Range("N276").Value = '헐, 반 밖에 안 남았네.'
# row 277
$ws.Range("A277").Value = 45570.64014188657
$ws.Range("B277").Value = 'kbi70722@gmail.com'
$ws.Range("C277").Value = '일본학과'
$ws.Range("D277").Value = 20191604
$ws.Range("E277").Value = '김병일'
$ws.Range("F277").Value = '민주 문자'
$ws.Range("G277").Value = '한글'
$ws.Range("H277").Value = '하나도 없다'
$ws.Range("I277").Value = 0.9
$ws.Range("J277").Value = '대한민국'
$ws.Range("K277").Value = '시간당 중위 임금이 60% 낮다'
$ws.Range("L277").Value = 'Red'
$ws.Range("M277").Value = '휴우, 그래도 반이나 남았네.'
# row 278
$ws.Range("A278").Value = 45570.64079939815
$ws.Range("B278").Value = 'quddus6378@gmail.com'
$ws.Range("C278").Value = '체육학과'
$ws.Range("D278").Value = 20227106
$ws.Range("E278").Value = '현병연'
$ws.Range("F278").Value = '민주 문자'
$ws.Range("G278").Value = '한글'
$ws.Range("H278").Value = '하나도 없다'
$ws.Range("I278").Value = 0.2
$ws.Range("J278").Value = '대한민국'
$ws.Range("K278").Value = '2배 정도 실직할 가능성이 높다'
$ws.Range("L278").Value = 'Black'
$ws.Range("N278").Value = '모름/기타'
# row 279
$ws.Range("A279").Value = 45570.6526321412
$ws.Range("B279").Value = '0524psu@gmail.com'
$ws.Range("C279").Value = '환경생명공학과'
$ws.Range("D279").Value = 20243712
$ws.Range("E279").Value = '박상언'
$ws.Range("F279").Value = '민주 문자'
$ws.Range("G279").Value = '한글'
$ws.Range("H279").Value = '하나도 없다'
$ws.Range("I279").Value = 0.8
$ws.Range("J279").Value = '대한민국'
$ws.Range("K279").Value = '시간당 중위 임금이 60% 낮다'
$ws.Range("L279").Value = 'Red'
$ws.Range("M279").Value = '휴우, 그래도 반이나 남았네.'
# row 280
$ws.Range("A280").Value = 45570.65527232639
$ws.Range("B280").Value = 'wlsqhwlsqh21@naver.com'
$ws.Range("C280").Value = '경영학부'
$ws.Range("D280").Value = 20242982
$ws.Range("E280").Value = '심진보'
$ws.Range("F280").Value = '민주 문자'
$ws.Range("G280").Value = '한글'
$ws.Range("H280").Value = '1개'
$ws.Range("I280").Value = 0.8
$ws.Range("J280").Value = '대한민국'
$ws.Range("K280").Value = '시간당 중위 임금이 60% 낮다'
$ws.Range("L280").Value = 'Black'
$ws.Range("N280").Value = '헐, 반 밖에 안 남았네.'
# row 281
$ws.Range("A281").Value = 45570.675174143515
$ws.Range("B281").Value = 'ella2005710@gmail.com'
$ws.Range("C281").Value = '간호학과'
$ws.Range("D281").Value = 20246222
$ws.Range("E281").Value = '김송이'
$ws.Range("F281").Value = '민주 문자'
$ws.Range("G281").Value = '한글'
$ws.Range("H281").Value = '3개'
$ws.Range("I281").Value = 0.8
$ws.Range("J281").Value = '대한민국'
$ws.Range("K281").Value = '시간당 중위 임금이 60% 낮다'
$ws.Range("L281").Value = 'Black'
$ws.Range("N281").Value = '헐, 반 밖에 안 남았네.'
# row 282
$ws.Range("A282").Value = 45570.67686912037
$ws.Range("B282").Value = 'dbfrhr02@naver.com'
$ws.Range("C282").Value = '경영학과'
$ws.Range("D282").Value = 20212971
$ws.Range("E282").Value = '심건휘'
$ws.Range("F282").Value = '민주 문자'
$ws.Range("G282").Value = '한글'
$ws.Range("H282").Value = '1개'
$ws.Range("I282").Value = 0.8
$ws.Range("J282").Value = '대한민국'
$ws.Range("K282").Value = '2배 정도 실직할 가능성이 높다'
$ws.Range("L282").Value = 'Red'
$ws.Range("M282").Value = '휴우, 그래도 반이나 남았네.'
# row 283
$ws.Range("A283").Value = 45570.68867167824
$ws.Range("B283").Value = '1kdcf@naver.com'
$ws.Range("C283").Value = '경영학과'
$ws.Range("D283").Value = 20192901
$ws.Range("E283").Value = '송준영'
$ws.Range("F283").Value = '민주 문자'
$ws.Range("G283").Value = '한글'
$ws.Range("H283").Value = '1개'
$ws.Range("I283").Value = 0.8
$ws.Range("J283").Value = '대한민국'
$ws.Range("K283").Value = '남들을 덜 신뢰한다'
$ws.Range("L283").Value = 'Red'
$ws.Range("M283").Value = '휴우, 그래도 반이나 남았네.'
# row 284
$ws.Range("A284").Value = 45570.69673954861
$ws.Range("B284").Value = 'hm703711@gmail.com'
$ws.Range("C284").Value = '언론방송융합미디어학과'
$ws.Range("D284").Value = 20202415
$ws.Range("E284").Value = '박현민'
$ws.Range("F284").Value = '민주 문자'
$ws.Range("G284").Value = '한글'
$ws.Range("H284").Value = '하나도 없다'
$ws.Range("I284").Value = 0.5
$ws.Range("J284").Value = '이탈리아'
$ws.Range("K284").Value = '남들을 덜 신뢰한다'
$ws.Range("L284").Value = 'Black'
$ws.Range("N284").Value = '헐, 반 밖에 안 남았네.'
# row 285
$ws.Range("A285").Value = 45570.70106546296
$ws.Range("B285").Value = 'gangjunu@naver.com'
$ws.Range("C285").Value = '금융재무학과'
$ws.Range("D285").Value = 20242901
$ws.Range("E285").Value = '강준우'
$ws.Range("F285").Value = '민주 문자'
$ws.Range("G285").Value = '한글'
$ws.Range("H285").Value = '하나도 없다'
$ws.Range("I285").Value = 0.9
$ws.Range("J285").Value = '대한민국'
$ws.Range("K285").Value = '남들을 덜 신뢰한다'
$ws.Range("L285").Value = 'Red'
$ws.Range("M285").Value = '휴우, 그래도 반이나 남았네.'
# row 286
$ws.Range("A286").Value = 45570.706889803245
$ws.Range("B286").Value = 'towp7563@gmail.com'
$ws.Range("C286").Value = '식품영양학과'
$ws.Range("D286").Value = 20243826
$ws.Range("E286").Value = '오승현'
$ws.Range("F286").Value = '민주 문자'
$ws.Range("G286").Value = '한글'
$ws.Range("H286").Value = '1개'
$ws.Range("I286").Value = 0.8
$ws.Range("J286").Value = '대한민국'
$ws.Range("K286").Value = '시간당 중위 임금이 60% 낮다'
$ws.Range("L286").Value = 'Black'
$ws.Range("N286").Value = '헐, 반 밖에 안 남았네.'
# row 287
$ws.Range("A287").Value = 45570.71745517361
$ws.Range("B287").Value = 'xodet0817@naver.com'
$ws.Range("C287").Value = '바이오메디컬학과'
$ws.Range("D287").Value = 20203616
$ws.Range("E287").Value = '문종윤'
$ws.Range("F287").Value = '민주 문자'
$ws.Range("G287").Value = '한글'
$ws.Range("H287").Value = '1개'
$ws.Range("I287").Value = 0.8
$ws.Range("J287").Value = '대한민국'
$ws.Range("K287").Value = '시간당 중위 임금이 60% 낮다'
$ws.Range("L287").Value = 'Black'
$ws.Range("N287").Value = '헐, 반 밖에 안 남았네.'
# row 288
$ws.Range("A288").Value = 45570.72358883102
$ws.Range("B288").Value = 'eojeongmin146@gmail.com'
$ws.Range("C288").Value = '법학과'
$ws.Range("D288").Value = 20242725
$ws.Range("E288").Value = '어정민'
$ws.Range("F288").Value = '민주 문자'
$ws.Range("G288").Value = '한글'
$ws.Range("H288").Value = '하나도 없다'
$ws.Range("I288").Value = 0.8
$ws.Range("J288").Value = '대한민국'
$ws.Range("K288").Value = '남들을 덜 신뢰한다'
$ws.Range("L288").Value = 'Red'
$ws.Range("M288").Value = '휴우, 그래도 반이나 남았네.'
# row 289
$ws.Range("A289").Value = 45570.72459525463
$ws.Range("B289").Value = 'seoeunchan5@gmail.com'
$ws.Range("C289").Value = '소프트웨어학과'
$ws.Range("D289").Value = 20245183
$ws.Range("E289").Value = '서은찬'
$ws.Range("F289").Value = '민주 문자'
$ws.Range("G289").Value = '한글'
$ws.Range("H289").Value = '2개'
$ws.Range("I289").Value = 0.8
$ws.Range("J289").Value = '대한민국'
$ws.Range("K289").Value = '시간당 중위 임금이 60% 낮다'
$ws.Range("L289").Value = 'Black'
$ws.Range("N289").Value = '헐, 반 밖에 안 남았네.'
# row 290
$ws.Range("A290").Value = 45570.72779075231
$ws.Range("B290").Value = 'jyb051128@gmail.com'
$ws.Range("C290").Value = '바이오메디컬학과'
$ws.Range("D290").Value = 20243646
$ws.Range("E290").Value = '전유빈'
$ws.Range("F290").Value = '민주 문자'
$ws.Range("G290").Value = '한자'
$ws.Range("H290").Value = '하나도 없다'
$ws.Range("I290").Value = 0.8
$ws.Range("J290").Value = '영국'
$ws.Range("K290").Value = '남들을 덜 신뢰한다'
$ws.Range("L290").Value = 'Black'
$ws.Range("N290").Value = '휴우, 그래도 반이나 남았네.'
# row 291
$ws.Range("A291").Value = 45570.72934377315
$ws.Range("B291").Value = 'rer220@naver.com'
$ws.Range("C291").Value = '콘텐츠IT'
$ws.Range("D291").Value = 20205124
$ws.Range("E291").Value = '김대명'
$ws.Range("F291").Value = '민주 문자'
$ws.Range("G291").Value = '한글'
$ws.Range("H291").Value = '하나도 없다'
$ws.Range("I291").Value = 0.8
$ws.Range("J291").Value = '대한민국'
$ws.Range("K291").Value = '시간당 중위 임금이 60% 낮다'
$ws.Range("L291").Value = 'Black'
$ws.Range("N291").Value = '헐, 반 밖에 안 남았네.'
# row 292
$ws.Range("A292").Value = 45570.73674439815
$ws.Range("B292").Value = 'gmlfkr6241@naver.com'
$ws.Range("C292").Value = '간호학과'
$ws.Range("D292").Value = 20246303
$ws.Range("E292").Value = '한희락'
$ws.Range("F292").Value = '민주 문자'
$ws.Range("G292").Value = '한글'
$ws.Range("H292").Value = '하나도 없다'
$ws.Range("I292").Value = 0.8
$ws.Range("J292").Value = '대한민국'
$ws.Range("K292").Value = '시간당 중위 임금이 60% 낮다'
$ws.Range("L292").Value = 'Black'
$ws.Range("N292").Value = '헐, 반 밖에 안 남았네.'
# row 293
$ws.Range("A293").Value = 45570.755309722226
$ws.Range("B293").Value = 'benjamin27@naver.com'
$ws.Range("C293").Value = '디지털미디어스쿨'
$ws.Range("D293").Value = 20212583
$ws.Range("E293").Value = '최재혁'
$ws.Range("F293").Value = '민주 문자'
$ws.Range("G293").Value = '한글'
$ws.Range("H293").Value = '2개'
$ws.Range("I293").Value = 0.8
$ws.Range("J293").Value = '영국'
$ws.Range("K293").Value = '시간당 중위 임금이 60% 낮다'
$ws.Range("L293").Value = 'Red'
$ws.Range("M293").Value = '휴우, 그래도 반이나 남았네.'

# 4) Grow the worksheet Table (ListObject) so it covers the new rows.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:N293"))
